$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text formatting so
# values such as "95.20" or "1.00" are not reinterpreted as numbers and
# lose their trailing zeros / thousands separators.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '44.563.60'
$ws.Range("E2").Value = '  +1.31%  '
$ws.Range("D3").Value = '2.242.14'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").Value = '306.51'
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").Value = '95.20'
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("D7").Value = '0.571'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = '0.521'
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("D10").Value = '34.93'
$ws.Range("E10").Value = '  +0.66%  '
$ws.Range("D11").Value = '0.0805'
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("D12").Value = '7.24'
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D14").Value = '2.274.77'
$ws.Range("E14").Value = '  +1.75%  '
$ws.Range("D15").Value = '0.834'
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").Value = '13.61'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").Value = '44.286.79'
$ws.Range("E17").Value = '  +0.96%  '
$ws.Range("D18").Value = '0.0₃0954'
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("D19").Value = '6.32'
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("D20").Value = '11.93'
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("D21").Value = '65.51'
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").Value = '237.15'
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").Value = '2.96'
$ws.Range("E23").Value = '  +1.20%  '
$ws.Range("D24").Value = '1.99'
$ws.Range("E24").Value = '  +1.33%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").Value = '2.26'
$ws.Range("E26").Value = '  +4.19%  '
$ws.Range("D27").Value = '9.81'
$ws.Range("E27").Value = '  -0.53%  '
$ws.Range("D28").Value = '37.59'
$ws.Range("E28").Value = '  -0.40%  '
$ws.Range("D29").Value = '5.97'
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = '19.93'
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("D31").Value = '152.44'
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").Value = '0.0798'
$ws.Range("E32").Value = '  -0.19%  '
$ws.Range("E33").Value = '  +2.05%  '
$ws.Range("D34").Value = '3.04'
$ws.Range("E34").Value = '  -6.41%  '
$ws.Range("D35").Value = '0.109'
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").Value = '1.84'
$ws.Range("E37").Value = '  +4.04%  '
$ws.Range("D38").Value = '14.91'
$ws.Range("E38").Value = '  -0.24%  '
$ws.Range("D39").Value = '3.38'
$ws.Range("E39").Value = '  +1.70%  '
$ws.Range("D40").Value = '3.77'
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("D41").Value = '0.0301'
$ws.Range("E41").Value = '  +1.26%  '
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("D43").Value = '1.806.46'
$ws.Range("E43").Value = '  +4.98%  '
$ws.Range("D44").Value = '1.69'
$ws.Range("E44").Value = '  +13.44%  '
$ws.Range("D45").Value = '0.192'
$ws.Range("E45").Value = '  +3.06%  '
$ws.Range("D46").Value = '78.75'
$ws.Range("E46").Value = '  -6.88%  '
$ws.Range("D47").Value = '70.59'
$ws.Range("E47").Value = '  +3.06%  '
$ws.Range("D48").Value = '98.87'
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("D49").Value = '4.90'
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("D51").Value = '54.41'
$ws.Range("E51").Value = '  +0.76%  '
